$wb = $excel.ActiveWorkbook

# --- "Projects" sheet: add column C "color" ---
$ws = $wb.Worksheets.Item("Projects")

# Header cell C1 should reuse the same header style as A1/B1.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C1").Value = "color"

# Data cells C2:C5 get the plain (unstyled) value.
$ws.Range("C2").Value = "#e0f7fa"
$ws.Range("C3").Value = "#e0f7fa"
$ws.Range("C4").Value = "#e0f7fa"
$ws.Range("C5").Value = "#e0f7fa"

# --- "Tasks" sheet: fill in D4 notes ---
$ws2 = $wb.Worksheets.Item("Tasks")
$ws2.Range("D4").Value = "<p>sdkj</p><p>alsk</p><ol><li>sd </li><li><br></li></ol>"
